# Error Calculations and Plots
# Applies the changes described by the diff:
#  - Removes two rows ("RM 232" and "SC 92") causing subsequent rows to
#    shift up (dimension A1:F35 -> A1:F33)
#  - Updates a handful of individual cell values (some filled in, some
#    cleared) to reflect the re-run imputation pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two removed rows. Delete the higher-numbered row first so the
# lower row index used for the second delete is still valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Individual cell updates on the resulting (shifted) sheet.
$ws.Range("F3").Value = 17.64
$ws.Range("D5").ClearContents()
$ws.Range("E6").Value = -5.7
$ws.Range("D8").Value = -13.9
$ws.Range("F10").ClearContents()
$ws.Range("E11").Value = -7.9
$ws.Range("F11").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("D14").Value = -13.1
$ws.Range("F16").Value = 17.34
$ws.Range("E17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("F20").Value = 17.73
$ws.Range("F24").ClearContents()
$ws.Range("E25").Value = -7.1
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("F28").ClearContents()
$ws.Range("E31").ClearContents()
$ws.Range("E32").ClearContents()
$ws.Range("F32").Value = 17.39
$ws.Range("C33").Value = 10.4
$ws.Range("F33").Value = 17.53
